# Update cryptos list (price + 1h volume change) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.026.05"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.237.79"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'305.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.88%  "
$ws.Range("E6").Value = "  -6.86%  "
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").Value = "'34.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.81%  "
$ws.Range("D11").Value = "'0.0807"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("E12").Value = "  -5.85%  "
$ws.Range("D13").Value = "'0.103"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "2.578.79"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "2.282.26"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "'0.826"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").Value = "'13.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.67%  "
$ws.Range("D18").Value = "43.888.66"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("E20").Value = "  -8.65%  "
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("D22").Value = "'64.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "'236.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'2.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.03%  "
$ws.Range("E25").Value = "  -7.87%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'9.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "'37.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("D29").Value = "'2.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").Value = "'20.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "'155.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").Value = "'3.28"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.69%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("E36").Value = "  -6.29%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("E38").Value = "  -11.08%  "
$ws.Range("D39").Value = "'15.32"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.02%  "
$ws.Range("E40").Value = "  -10.16%  "
$ws.Range("E41").Value = "  -10.04%  "
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "1.735.66"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "'85.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("D47").Value = "'99.89"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").Value = "'4.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.99%  "
$ws.Range("D49").Value = "'69.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.32%  "
$ws.Range("D50").Value = "'8.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").Value = "'54.22"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.58%  "
